$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value() -eq "stimuli/blank.JPG") {
        $cell.Value() = "stimuli/blank.jpg"
    }
}
